$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select cell A2 (reflected as the sheet's active selection)
$ws.Range("A2").Select()

# Fix priors and parameters for stage distribution
$ws.Range("E6").Value = -0.05
$ws.Range("F6").Value = 0.03

$ws.Range("E7").Value = -0.03
$ws.Range("F7").Value = 0.03

$ws.Range("E8").Value = -0.03
$ws.Range("F8").Value = 0.03
